$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colD = New-Object "object[,]" 85,1
$colD[0,0] = 45223
$colD[1,0] = 45223
$colD[2,0] = 44468
$colD[3,0] = 44487
$colD[4,0] = 44487
$colD[5,0] = 44641
$colD[6,0] = 44714
$colD[7,0] = 44714
$colD[8,0] = 44714
$colD[9,0] = 44778
$colD[10,0] = 44778
$colD[11,0] = 44754
$colD[12,0] = 44754
$colD[13,0] = 44754
$colD[14,0] = 44462
$colD[15,0] = 44462
$colD[16,0] = 44462
$colD[17,0] = 44490
$colD[18,0] = 44490
$colD[19,0] = 45215
$colD[20,0] = 44238
$colD[21,0] = 44238
$colD[22,0] = 44973
$colD[23,0] = 44973
$colD[24,0] = 44874
$colD[25,0] = 44874
$colD[26,0] = 44245
$colD[27,0] = 44245
$colD[28,0] = 44481
$colD[29,0] = 44481
$colD[30,0] = 44999
$colD[31,0] = 44999
$colD[32,0] = 44294
$colD[33,0] = 44294
$colD[34,0] = 44558
$colD[35,0] = 44634
$colD[36,0] = 44634
$colD[37,0] = 44174
$colD[38,0] = 44608
$colD[39,0] = 44550
$colD[40,0] = 44775
$colD[41,0] = 44775
$colD[42,0] = 44859
$colD[43,0] = 44859
$colD[44,0] = 44859
$colD[45,0] = 44910
$colD[46,0] = 44910
$colD[47,0] = 45063
$colD[48,0] = 44365
$colD[49,0] = 44365
$colD[50,0] = 44365
$colD[51,0] = 45069
$colD[52,0] = 44711
$colD[53,0] = 44382
$colD[54,0] = 44382
$colD[55,0] = 44382
$colD[56,0] = 44795
$colD[57,0] = 45173
$colD[58,0] = 44883
$colD[59,0] = 44883
$colD[60,0] = 44673
$colD[61,0] = 44818
$colD[62,0] = 45140
$colD[63,0] = 44649
$colD[64,0] = 44649
$colD[65,0] = 44341
$colD[66,0] = 44341
$colD[67,0] = 44341
$colD[68,0] = 44868
$colD[69,0] = 45216
$colD[70,0] = 45216
$colD[71,0] = 45216
$colD[72,0] = 44460
$colD[73,0] = 45114
$colD[74,0] = 45114
$colD[75,0] = 45012
$colD[76,0] = 44232
$colD[77,0] = 44232
$colD[78,0] = 44613
$colD[79,0] = 44725
$colD[80,0] = 44893
$colD[81,0] = 44392
$colD[82,0] = 44392
$colD[83,0] = 44565
$colD[84,0] = 44544
$ws.Range("D935:D1019").Value = $colD

$blockHQ = New-Object "object[,]" 85,10
$blockHQ[0,0] = 'Escarola'
$blockHQ[0,1] = 'Primera'
$blockHQ[0,2] = 500
$blockHQ[0,3] = 19000
$blockHQ[0,4] = 19000
$blockHQ[0,5] = 19000
$blockHQ[0,6] = '$/caja 15 unidades'
$blockHQ[0,7] = 'Región de Coquimbo'
$blockHQ[0,8] = 1267
$blockHQ[0,9] = 15
$blockHQ[1,0] = 'Marina'
$blockHQ[1,1] = 'Primera'
$blockHQ[1,2] = 300
$blockHQ[1,3] = 12000
$blockHQ[1,4] = 12000
$blockHQ[1,5] = 12000
$blockHQ[1,6] = '$/caja 15 unidades'
$blockHQ[1,7] = 'Región Metropolitana'
$blockHQ[1,8] = 800
$blockHQ[1,9] = 15
$blockHQ[2,0] = 'Escarola'
$blockHQ[2,1] = 'Primera'
$blockHQ[2,2] = 150
$blockHQ[2,3] = 12000
$blockHQ[2,4] = 12000
$blockHQ[2,5] = 12000
$blockHQ[2,6] = '$/caja 15 unidades'
$blockHQ[2,7] = 'Región de Coquimbo'
$blockHQ[2,8] = 800
$blockHQ[2,9] = 15
$blockHQ[3,0] = 'Conconina(o)'
$blockHQ[3,1] = 'Primera'
$blockHQ[3,2] = 120
$blockHQ[3,3] = 8000
$blockHQ[3,4] = 8000
$blockHQ[3,5] = 8000
$blockHQ[3,6] = '$/caja 10 unidades'
$blockHQ[3,7] = 'Región del Maule'
$blockHQ[3,8] = 800
$blockHQ[3,9] = 10
$blockHQ[4,0] = 'Escarola'
$blockHQ[4,1] = 'Primera'
$blockHQ[4,2] = 250
$blockHQ[4,3] = 10000
$blockHQ[4,4] = 11000
$blockHQ[4,5] = 10600
$blockHQ[4,6] = '$/caja 15 unidades'
$blockHQ[4,7] = 'Región de Coquimbo'
$blockHQ[4,8] = 707
$blockHQ[4,9] = 15
$blockHQ[5,0] = 'Escarola'
$blockHQ[5,1] = 'Primera'
$blockHQ[5,2] = 300
$blockHQ[5,3] = 12000
$blockHQ[5,4] = 14000
$blockHQ[5,5] = 13000
$blockHQ[5,6] = '$/caja 15 unidades'
$blockHQ[5,7] = 'Región del Maule'
$blockHQ[5,8] = 867
$blockHQ[5,9] = 15
$blockHQ[6,0] = 'Conconina(o)'
$blockHQ[6,1] = 'Primera'
$blockHQ[6,2] = 150
$blockHQ[6,3] = 12000
$blockHQ[6,4] = 12000
$blockHQ[6,5] = 12000
$blockHQ[6,6] = '$/caja 10 unidades'
$blockHQ[6,7] = 'Región Metropolitana'
$blockHQ[6,8] = 1200
$blockHQ[6,9] = 10
$blockHQ[7,0] = 'Escarola'
$blockHQ[7,1] = 'Primera'
$blockHQ[7,2] = 150
$blockHQ[7,3] = 12000
$blockHQ[7,4] = 12000
$blockHQ[7,5] = 12000
$blockHQ[7,6] = '$/caja 15 unidades'
$blockHQ[7,7] = 'Región de Coquimbo'
$blockHQ[7,8] = 800
$blockHQ[7,9] = 15
$blockHQ[8,0] = 'Escarola'
$blockHQ[8,1] = 'Segunda'
$blockHQ[8,2] = 150
$blockHQ[8,3] = 10000
$blockHQ[8,4] = 10000
$blockHQ[8,5] = 10000
$blockHQ[8,6] = '$/caja 18 unidades'
$blockHQ[8,7] = 'Región de Coquimbo'
$blockHQ[8,8] = 556
$blockHQ[8,9] = 18
$blockHQ[9,0] = 'Escarola'
$blockHQ[9,1] = 'Primera'
$blockHQ[9,2] = 300
$blockHQ[9,3] = 15000
$blockHQ[9,4] = 15000
$blockHQ[9,5] = 15000
$blockHQ[9,6] = '$/caja 15 unidades'
$blockHQ[9,7] = 'Región de Coquimbo'
$blockHQ[9,8] = 1000
$blockHQ[9,9] = 15
$blockHQ[10,0] = 'Escarola'
$blockHQ[10,1] = 'Segunda'
$blockHQ[10,2] = 300
$blockHQ[10,3] = 12000
$blockHQ[10,4] = 12000
$blockHQ[10,5] = 12000
$blockHQ[10,6] = '$/caja 18 unidades'
$blockHQ[10,7] = 'Región de Coquimbo'
$blockHQ[10,8] = 667
$blockHQ[10,9] = 18
$blockHQ[11,0] = 'Escarola'
$blockHQ[11,1] = 'Primera'
$blockHQ[11,2] = 300
$blockHQ[11,3] = 15000
$blockHQ[11,4] = 15000
$blockHQ[11,5] = 15000
$blockHQ[11,6] = '$/caja 15 unidades'
$blockHQ[11,7] = 'Región de Coquimbo'
$blockHQ[11,8] = 1000
$blockHQ[11,9] = 15
$blockHQ[12,0] = 'Escarola'
$blockHQ[12,1] = 'Segunda'
$blockHQ[12,2] = 400
$blockHQ[12,3] = 12500
$blockHQ[12,4] = 13000
$blockHQ[12,5] = 12750
$blockHQ[12,6] = '$/caja 18 unidades'
$blockHQ[12,7] = 'Región de Coquimbo'
$blockHQ[12,8] = 708
$blockHQ[12,9] = 18
$blockHQ[13,0] = 'Marina'
$blockHQ[13,1] = 'Primera'
$blockHQ[13,2] = 350
$blockHQ[13,3] = 9000
$blockHQ[13,4] = 9500
$blockHQ[13,5] = 9214
$blockHQ[13,6] = '$/caja 15 unidades'
$blockHQ[13,7] = 'Región Metropolitana'
$blockHQ[13,8] = 614
$blockHQ[13,9] = 15
$blockHQ[14,0] = 'Conconina(o)'
$blockHQ[14,1] = 'Primera'
$blockHQ[14,2] = 80
$blockHQ[14,3] = 9000
$blockHQ[14,4] = 9000
$blockHQ[14,5] = 9000
$blockHQ[14,6] = '$/caja 10 unidades'
$blockHQ[14,7] = 'Región Metropolitana'
$blockHQ[14,8] = 900
$blockHQ[14,9] = 10
$blockHQ[15,0] = 'Escarola'
$blockHQ[15,1] = 'Primera'
$blockHQ[15,2] = 250
$blockHQ[15,3] = 12000
$blockHQ[15,4] = 12000
$blockHQ[15,5] = 12000
$blockHQ[15,6] = '$/caja 15 unidades'
$blockHQ[15,7] = 'Región de Coquimbo'
$blockHQ[15,8] = 800
$blockHQ[15,9] = 15
$blockHQ[16,0] = 'Marina'
$blockHQ[16,1] = 'Primera'
$blockHQ[16,2] = 150
$blockHQ[16,3] = 9000
$blockHQ[16,4] = 9000
$blockHQ[16,5] = 9000
$blockHQ[16,6] = '$/caja 15 unidades'
$blockHQ[16,7] = 'Región Metropolitana'
$blockHQ[16,8] = 600
$blockHQ[16,9] = 15
$blockHQ[17,0] = 'Conconina(o)'
$blockHQ[17,1] = 'Primera'
$blockHQ[17,2] = 100
$blockHQ[17,3] = 8000
$blockHQ[17,4] = 9000
$blockHQ[17,5] = 8500
$blockHQ[17,6] = '$/caja 10 unidades'
$blockHQ[17,7] = 'Región Metropolitana'
$blockHQ[17,8] = 850
$blockHQ[17,9] = 10
$blockHQ[18,0] = 'Escarola'
$blockHQ[18,1] = 'Primera'
$blockHQ[18,2] = 250
$blockHQ[18,3] = 10000
$blockHQ[18,4] = 11000
$blockHQ[18,5] = 10600
$blockHQ[18,6] = '$/caja 15 unidades'
$blockHQ[18,7] = 'Región de Coquimbo'
$blockHQ[18,8] = 707
$blockHQ[18,9] = 15
$blockHQ[19,0] = 'Escarola'
$blockHQ[19,1] = 'Primera'
$blockHQ[19,2] = 120
$blockHQ[19,3] = 18000
$blockHQ[19,4] = 18000
$blockHQ[19,5] = 18000
$blockHQ[19,6] = '$/caja 15 unidades'
$blockHQ[19,7] = 'Región de Coquimbo'
$blockHQ[19,8] = 1200
$blockHQ[19,9] = 15
$blockHQ[20,0] = 'Conconina(o)'
$blockHQ[20,1] = 'Primera'
$blockHQ[20,2] = 100
$blockHQ[20,3] = 8000
$blockHQ[20,4] = 9000
$blockHQ[20,5] = 8500
$blockHQ[20,6] = '$/caja 10 unidades'
$blockHQ[20,7] = 'Región Metropolitana'
$blockHQ[20,8] = 850
$blockHQ[20,9] = 10
$blockHQ[21,0] = 'Escarola'
$blockHQ[21,1] = 'Primera'
$blockHQ[21,2] = 200
$blockHQ[21,3] = 8000
$blockHQ[21,4] = 9000
$blockHQ[21,5] = 8500
$blockHQ[21,6] = '$/caja 15 unidades'
$blockHQ[21,7] = 'Región de Coquimbo'
$blockHQ[21,8] = 567
$blockHQ[21,9] = 15
$blockHQ[22,0] = 'Conconina(o)'
$blockHQ[22,1] = 'Primera'
$blockHQ[22,2] = 80
$blockHQ[22,3] = 12000
$blockHQ[22,4] = 12000
$blockHQ[22,5] = 12000
$blockHQ[22,6] = '$/caja 10 unidades'
$blockHQ[22,7] = 'Región Metropolitana'
$blockHQ[22,8] = 1200
$blockHQ[22,9] = 10
$blockHQ[23,0] = 'Escarola'
$blockHQ[23,1] = 'Primera'
$blockHQ[23,2] = 200
$blockHQ[23,3] = 15000
$blockHQ[23,4] = 15000
$blockHQ[23,5] = 15000
$blockHQ[23,6] = '$/caja 15 unidades'
$blockHQ[23,7] = 'Región de Coquimbo'
$blockHQ[23,8] = 1000
$blockHQ[23,9] = 15
$blockHQ[24,0] = 'Escarola'
$blockHQ[24,1] = 'Primera'
$blockHQ[24,2] = 100
$blockHQ[24,3] = 9000
$blockHQ[24,4] = 9000
$blockHQ[24,5] = 9000
$blockHQ[24,6] = '$/caja 15 unidades'
$blockHQ[24,7] = 'Región Metropolitana'
$blockHQ[24,8] = 600
$blockHQ[24,9] = 15
$blockHQ[25,0] = 'Escarola'
$blockHQ[25,1] = 'Segunda'
$blockHQ[25,2] = 100
$blockHQ[25,3] = 7000
$blockHQ[25,4] = 7000
$blockHQ[25,5] = 7000
$blockHQ[25,6] = '$/caja 18 unidades'
$blockHQ[25,7] = 'Región Metropolitana'
$blockHQ[25,8] = 389
$blockHQ[25,9] = 18
$blockHQ[26,0] = 'Conconina(o)'
$blockHQ[26,1] = 'Primera'
$blockHQ[26,2] = 60
$blockHQ[26,3] = 8000
$blockHQ[26,4] = 8000
$blockHQ[26,5] = 8000
$blockHQ[26,6] = '$/caja 10 unidades'
$blockHQ[26,7] = 'Región Metropolitana'
$blockHQ[26,8] = 800
$blockHQ[26,9] = 10
$blockHQ[27,0] = 'Escarola'
$blockHQ[27,1] = 'Primera'
$blockHQ[27,2] = 150
$blockHQ[27,3] = 10000
$blockHQ[27,4] = 10000
$blockHQ[27,5] = 10000
$blockHQ[27,6] = '$/caja 15 unidades'
$blockHQ[27,7] = 'Región de Coquimbo'
$blockHQ[27,8] = 667
$blockHQ[27,9] = 15
$blockHQ[28,0] = 'Escarola'
$blockHQ[28,1] = 'Primera'
$blockHQ[28,2] = 450
$blockHQ[28,3] = 9000
$blockHQ[28,4] = 9000
$blockHQ[28,5] = 9000
$blockHQ[28,6] = '$/caja 15 unidades'
$blockHQ[28,7] = 'Región Metropolitana'
$blockHQ[28,8] = 600
$blockHQ[28,9] = 15
$blockHQ[29,0] = 'Marina'
$blockHQ[29,1] = 'Segunda'
$blockHQ[29,2] = 250
$blockHQ[29,3] = 7000
$blockHQ[29,4] = 7000
$blockHQ[29,5] = 7000
$blockHQ[29,6] = '$/caja 18 unidades'
$blockHQ[29,7] = 'Región Metropolitana'
$blockHQ[29,8] = 389
$blockHQ[29,9] = 18
$blockHQ[30,0] = 'Escarola'
$blockHQ[30,1] = 'Primera'
$blockHQ[30,2] = 300
$blockHQ[30,3] = 12000
$blockHQ[30,4] = 12000
$blockHQ[30,5] = 12000
$blockHQ[30,6] = '$/caja 15 unidades'
$blockHQ[30,7] = 'Región de Coquimbo'
$blockHQ[30,8] = 800
$blockHQ[30,9] = 15
$blockHQ[31,0] = 'Escarola'
$blockHQ[31,1] = 'Segunda'
$blockHQ[31,2] = 250
$blockHQ[31,3] = 10500
$blockHQ[31,4] = 10500
$blockHQ[31,5] = 10500
$blockHQ[31,6] = '$/caja 18 unidades'
$blockHQ[31,7] = 'Región de Coquimbo'
$blockHQ[31,8] = 583
$blockHQ[31,9] = 18
$blockHQ[32,0] = 'Conconina(o)'
$blockHQ[32,1] = 'Primera'
$blockHQ[32,2] = 70
$blockHQ[32,3] = 11000
$blockHQ[32,4] = 11000
$blockHQ[32,5] = 11000
$blockHQ[32,6] = '$/caja 10 unidades'
$blockHQ[32,7] = 'Región Metropolitana'
$blockHQ[32,8] = 1100
$blockHQ[32,9] = 10
$blockHQ[33,0] = 'Escarola'
$blockHQ[33,1] = 'Primera'
$blockHQ[33,2] = 250
$blockHQ[33,3] = 11000
$blockHQ[33,4] = 11000
$blockHQ[33,5] = 11000
$blockHQ[33,6] = '$/caja 15 unidades'
$blockHQ[33,7] = 'Región de Coquimbo'
$blockHQ[33,8] = 733
$blockHQ[33,9] = 15
$blockHQ[34,0] = 'Escarola'
$blockHQ[34,1] = 'Primera'
$blockHQ[34,2] = 300
$blockHQ[34,3] = 10000
$blockHQ[34,4] = 10500
$blockHQ[34,5] = 10250
$blockHQ[34,6] = '$/caja 15 unidades'
$blockHQ[34,7] = 'Región Metropolitana'
$blockHQ[34,8] = 683
$blockHQ[34,9] = 15
$blockHQ[35,0] = 'Escarola'
$blockHQ[35,1] = 'Primera'
$blockHQ[35,2] = 150
$blockHQ[35,3] = 12000
$blockHQ[35,4] = 12000
$blockHQ[35,5] = 12000
$blockHQ[35,6] = '$/caja 15 unidades'
$blockHQ[35,7] = 'Región del Maule'
$blockHQ[35,8] = 800
$blockHQ[35,9] = 15
$blockHQ[36,0] = 'Escarola'
$blockHQ[36,1] = 'Segunda'
$blockHQ[36,2] = 100
$blockHQ[36,3] = 9000
$blockHQ[36,4] = 9000
$blockHQ[36,5] = 9000
$blockHQ[36,6] = '$/caja 18 unidades'
$blockHQ[36,7] = 'Región de La Araucanía'
$blockHQ[36,8] = 500
$blockHQ[36,9] = 18
$blockHQ[37,0] = 'Escarola'
$blockHQ[37,1] = 'Primera'
$blockHQ[37,2] = 27
$blockHQ[37,3] = 7000
$blockHQ[37,4] = 7000
$blockHQ[37,5] = 7000
$blockHQ[37,6] = '$/caja 15 unidades'
$blockHQ[37,7] = 'Región de La Araucanía'
$blockHQ[37,8] = 467
$blockHQ[37,9] = 15
$blockHQ[38,0] = 'Escarola'
$blockHQ[38,1] = 'Primera'
$blockHQ[38,2] = 100
$blockHQ[38,3] = 12000
$blockHQ[38,4] = 12000
$blockHQ[38,5] = 12000
$blockHQ[38,6] = '$/caja 15 unidades'
$blockHQ[38,7] = 'Región de Coquimbo'
$blockHQ[38,8] = 800
$blockHQ[38,9] = 15
$blockHQ[39,0] = 'Escarola'
$blockHQ[39,1] = 'Primera'
$blockHQ[39,2] = 200
$blockHQ[39,3] = 10000
$blockHQ[39,4] = 10000
$blockHQ[39,5] = 10000
$blockHQ[39,6] = '$/caja 15 unidades'
$blockHQ[39,7] = 'Región Metropolitana'
$blockHQ[39,8] = 667
$blockHQ[39,9] = 15
$blockHQ[40,0] = 'Escarola'
$blockHQ[40,1] = 'Primera'
$blockHQ[40,2] = 250
$blockHQ[40,3] = 15000
$blockHQ[40,4] = 15000
$blockHQ[40,5] = 15000
$blockHQ[40,6] = '$/caja 15 unidades'
$blockHQ[40,7] = 'Región de Coquimbo'
$blockHQ[40,8] = 1000
$blockHQ[40,9] = 15
$blockHQ[41,0] = 'Escarola'
$blockHQ[41,1] = 'Segunda'
$blockHQ[41,2] = 300
$blockHQ[41,3] = 12000
$blockHQ[41,4] = 13000
$blockHQ[41,5] = 12500
$blockHQ[41,6] = '$/caja 18 unidades'
$blockHQ[41,7] = 'Región de Coquimbo'
$blockHQ[41,8] = 694
$blockHQ[41,9] = 18
$blockHQ[42,0] = 'Escarola'
$blockHQ[42,1] = 'Primera'
$blockHQ[42,2] = 300
$blockHQ[42,3] = 10000
$blockHQ[42,4] = 10000
$blockHQ[42,5] = 10000
$blockHQ[42,6] = '$/caja 15 unidades'
$blockHQ[42,7] = 'Región de Coquimbo'
$blockHQ[42,8] = 667
$blockHQ[42,9] = 15
$blockHQ[43,0] = 'Escarola'
$blockHQ[43,1] = 'Segunda'
$blockHQ[43,2] = 300
$blockHQ[43,3] = 8500
$blockHQ[43,4] = 8500
$blockHQ[43,5] = 8500
$blockHQ[43,6] = '$/caja 18 unidades'
$blockHQ[43,7] = 'Región de Coquimbo'
$blockHQ[43,8] = 472
$blockHQ[43,9] = 18
$blockHQ[44,0] = 'Marina'
$blockHQ[44,1] = 'Segunda'
$blockHQ[44,2] = 300
$blockHQ[44,3] = 6500
$blockHQ[44,4] = 6500
$blockHQ[44,5] = 6500
$blockHQ[44,6] = '$/caja 18 unidades'
$blockHQ[44,7] = 'Región Metropolitana'
$blockHQ[44,8] = 361
$blockHQ[44,9] = 18
$blockHQ[45,0] = 'Escarola'
$blockHQ[45,1] = 'Primera'
$blockHQ[45,2] = 200
$blockHQ[45,3] = 12000
$blockHQ[45,4] = 13000
$blockHQ[45,5] = 12500
$blockHQ[45,6] = '$/caja 15 unidades'
$blockHQ[45,7] = 'Región de Coquimbo'
$blockHQ[45,8] = 833
$blockHQ[45,9] = 15
$blockHQ[46,0] = 'Escarola'
$blockHQ[46,1] = 'Segunda'
$blockHQ[46,2] = 100
$blockHQ[46,3] = 11000
$blockHQ[46,4] = 11000
$blockHQ[46,5] = 11000
$blockHQ[46,6] = '$/caja 18 unidades'
$blockHQ[46,7] = 'Región de Coquimbo'
$blockHQ[46,8] = 611
$blockHQ[46,9] = 18
$blockHQ[47,0] = 'Escarola'
$blockHQ[47,1] = 'Primera'
$blockHQ[47,2] = 120
$blockHQ[47,3] = 11000
$blockHQ[47,4] = 11000
$blockHQ[47,5] = 11000
$blockHQ[47,6] = '$/caja 15 unidades'
$blockHQ[47,7] = 'Región de Coquimbo'
$blockHQ[47,8] = 733
$blockHQ[47,9] = 15
$blockHQ[48,0] = 'Conconina(o)'
$blockHQ[48,1] = 'Primera'
$blockHQ[48,2] = 100
$blockHQ[48,3] = 11000
$blockHQ[48,4] = 11000
$blockHQ[48,5] = 11000
$blockHQ[48,6] = '$/caja 10 unidades'
$blockHQ[48,7] = 'Región Metropolitana'
$blockHQ[48,8] = 1100
$blockHQ[48,9] = 10
$blockHQ[49,0] = 'Escarola'
$blockHQ[49,1] = 'Primera'
$blockHQ[49,2] = 400
$blockHQ[49,3] = 12000
$blockHQ[49,4] = 12000
$blockHQ[49,5] = 12000
$blockHQ[49,6] = '$/caja 15 unidades'
$blockHQ[49,7] = 'Región de Coquimbo'
$blockHQ[49,8] = 800
$blockHQ[49,9] = 15
$blockHQ[50,0] = 'Marina'
$blockHQ[50,1] = 'Segunda'
$blockHQ[50,2] = 200
$blockHQ[50,3] = 8500
$blockHQ[50,4] = 8500
$blockHQ[50,5] = 8500
$blockHQ[50,6] = '$/caja 18 unidades'
$blockHQ[50,7] = 'Región Metropolitana'
$blockHQ[50,8] = 472
$blockHQ[50,9] = 18
$blockHQ[51,0] = 'Escarola'
$blockHQ[51,1] = 'Primera'
$blockHQ[51,2] = 600
$blockHQ[51,3] = 9500
$blockHQ[51,4] = 10500
$blockHQ[51,5] = 10000
$blockHQ[51,6] = '$/caja 15 unidades'
$blockHQ[51,7] = 'Región de Coquimbo'
$blockHQ[51,8] = 667
$blockHQ[51,9] = 15
$blockHQ[52,0] = 'Escarola'
$blockHQ[52,1] = 'Primera'
$blockHQ[52,2] = 250
$blockHQ[52,3] = 12000
$blockHQ[52,4] = 12000
$blockHQ[52,5] = 12000
$blockHQ[52,6] = '$/caja 15 unidades'
$blockHQ[52,7] = 'Región de Coquimbo'
$blockHQ[52,8] = 800
$blockHQ[52,9] = 15
$blockHQ[53,0] = 'Conconina(o)'
$blockHQ[53,1] = 'Primera'
$blockHQ[53,2] = 70
$blockHQ[53,3] = 11000
$blockHQ[53,4] = 11000
$blockHQ[53,5] = 11000
$blockHQ[53,6] = '$/caja 10 unidades'
$blockHQ[53,7] = 'Región Metropolitana'
$blockHQ[53,8] = 1100
$blockHQ[53,9] = 10
$blockHQ[54,0] = 'Escarola'
$blockHQ[54,1] = 'Primera'
$blockHQ[54,2] = 200
$blockHQ[54,3] = 12000
$blockHQ[54,4] = 12000
$blockHQ[54,5] = 12000
$blockHQ[54,6] = '$/caja 15 unidades'
$blockHQ[54,7] = 'Región de Coquimbo'
$blockHQ[54,8] = 800
$blockHQ[54,9] = 15
$blockHQ[55,0] = 'Marina'
$blockHQ[55,1] = 'Segunda'
$blockHQ[55,2] = 90
$blockHQ[55,3] = 8500
$blockHQ[55,4] = 8500
$blockHQ[55,5] = 8500
$blockHQ[55,6] = '$/caja 18 unidades'
$blockHQ[55,7] = 'Región Metropolitana'
$blockHQ[55,8] = 472
$blockHQ[55,9] = 18
$blockHQ[56,0] = 'Escarola'
$blockHQ[56,1] = 'Primera'
$blockHQ[56,2] = 300
$blockHQ[56,3] = 14000
$blockHQ[56,4] = 15000
$blockHQ[56,5] = 14500
$blockHQ[56,6] = '$/caja 15 unidades'
$blockHQ[56,7] = 'Región de Coquimbo'
$blockHQ[56,8] = 967
$blockHQ[56,9] = 15
$blockHQ[57,0] = 'Escarola'
$blockHQ[57,1] = 'Primera'
$blockHQ[57,2] = 250
$blockHQ[57,3] = 12000
$blockHQ[57,4] = 12000
$blockHQ[57,5] = 12000
$blockHQ[57,6] = '$/caja 15 unidades'
$blockHQ[57,7] = 'Región de Coquimbo'
$blockHQ[57,8] = 800
$blockHQ[57,9] = 15
$blockHQ[58,0] = 'Escarola'
$blockHQ[58,1] = 'Primera'
$blockHQ[58,2] = 300
$blockHQ[58,3] = 8000
$blockHQ[58,4] = 8000
$blockHQ[58,5] = 8000
$blockHQ[58,6] = '$/caja 15 unidades'
$blockHQ[58,7] = 'Región de Coquimbo'
$blockHQ[58,8] = 533
$blockHQ[58,9] = 15
$blockHQ[59,0] = 'Escarola'
$blockHQ[59,1] = 'Segunda'
$blockHQ[59,2] = 300
$blockHQ[59,3] = 7000
$blockHQ[59,4] = 7000
$blockHQ[59,5] = 7000
$blockHQ[59,6] = '$/caja 18 unidades'
$blockHQ[59,7] = 'Región de Coquimbo'
$blockHQ[59,8] = 389
$blockHQ[59,9] = 18
$blockHQ[60,0] = 'Escarola'
$blockHQ[60,1] = 'Primera'
$blockHQ[60,2] = 300
$blockHQ[60,3] = 14000
$blockHQ[60,4] = 15000
$blockHQ[60,5] = 14500
$blockHQ[60,6] = '$/caja 15 unidades'
$blockHQ[60,7] = 'Región de Coquimbo'
$blockHQ[60,8] = 967
$blockHQ[60,9] = 15
$blockHQ[61,0] = 'Escarola'
$blockHQ[61,1] = 'Segunda'
$blockHQ[61,2] = 120
$blockHQ[61,3] = 12500
$blockHQ[61,4] = 12500
$blockHQ[61,5] = 12500
$blockHQ[61,6] = '$/caja 18 unidades'
$blockHQ[61,7] = 'Región de Coquimbo'
$blockHQ[61,8] = 694
$blockHQ[61,9] = 18
$blockHQ[62,0] = 'Escarola'
$blockHQ[62,1] = 'Primera'
$blockHQ[62,2] = 120
$blockHQ[62,3] = 12000
$blockHQ[62,4] = 12000
$blockHQ[62,5] = 12000
$blockHQ[62,6] = '$/caja 15 unidades'
$blockHQ[62,7] = 'Región de Coquimbo'
$blockHQ[62,8] = 800
$blockHQ[62,9] = 15
$blockHQ[63,0] = 'Escarola'
$blockHQ[63,1] = 'Segunda'
$blockHQ[63,2] = 400
$blockHQ[63,3] = 10000
$blockHQ[63,4] = 11000
$blockHQ[63,5] = 10500
$blockHQ[63,6] = '$/caja 18 unidades'
$blockHQ[63,7] = 'Región de Coquimbo'
$blockHQ[63,8] = 583
$blockHQ[63,9] = 18
$blockHQ[64,0] = 'Marina'
$blockHQ[64,1] = 'Segunda'
$blockHQ[64,2] = 250
$blockHQ[64,3] = 8000
$blockHQ[64,4] = 8000
$blockHQ[64,5] = 8000
$blockHQ[64,6] = '$/caja 18 unidades'
$blockHQ[64,7] = 'Región Metropolitana'
$blockHQ[64,8] = 444
$blockHQ[64,9] = 18
$blockHQ[65,0] = 'Conconina(o)'
$blockHQ[65,1] = 'Segunda'
$blockHQ[65,2] = 130
$blockHQ[65,3] = 7500
$blockHQ[65,4] = 7500
$blockHQ[65,5] = 7500
$blockHQ[65,6] = '$/caja 12 unidades'
$blockHQ[65,7] = 'Región Metropolitana'
$blockHQ[65,8] = 625
$blockHQ[65,9] = 12
$blockHQ[66,0] = 'Escarola'
$blockHQ[66,1] = 'Primera'
$blockHQ[66,2] = 300
$blockHQ[66,3] = 12000
$blockHQ[66,4] = 12000
$blockHQ[66,5] = 12000
$blockHQ[66,6] = '$/caja 15 unidades'
$blockHQ[66,7] = 'Región de Coquimbo'
$blockHQ[66,8] = 800
$blockHQ[66,9] = 15
$blockHQ[67,0] = 'Marina'
$blockHQ[67,1] = 'Segunda'
$blockHQ[67,2] = 250
$blockHQ[67,3] = 8500
$blockHQ[67,4] = 8500
$blockHQ[67,5] = 8500
$blockHQ[67,6] = '$/caja 18 unidades'
$blockHQ[67,7] = 'Región Metropolitana'
$blockHQ[67,8] = 472
$blockHQ[67,9] = 18
$blockHQ[68,0] = 'Escarola'
$blockHQ[68,1] = 'Primera'
$blockHQ[68,2] = 300
$blockHQ[68,3] = 8500
$blockHQ[68,4] = 8500
$blockHQ[68,5] = 8500
$blockHQ[68,6] = '$/caja 15 unidades'
$blockHQ[68,7] = 'Región de Coquimbo'
$blockHQ[68,8] = 567
$blockHQ[68,9] = 15
$blockHQ[69,0] = 'Conconina(o)'
$blockHQ[69,1] = 'Primera'
$blockHQ[69,2] = 150
$blockHQ[69,3] = 13000
$blockHQ[69,4] = 13000
$blockHQ[69,5] = 13000
$blockHQ[69,6] = '$/caja 10 unidades'
$blockHQ[69,7] = 'Región Metropolitana'
$blockHQ[69,8] = 1300
$blockHQ[69,9] = 10
$blockHQ[70,0] = 'Escarola'
$blockHQ[70,1] = 'Primera'
$blockHQ[70,2] = 500
$blockHQ[70,3] = 19000
$blockHQ[70,4] = 19000
$blockHQ[70,5] = 19000
$blockHQ[70,6] = '$/caja 15 unidades'
$blockHQ[70,7] = 'Región de Coquimbo'
$blockHQ[70,8] = 1267
$blockHQ[70,9] = 15
$blockHQ[71,0] = 'Marina'
$blockHQ[71,1] = 'Primera'
$blockHQ[71,2] = 300
$blockHQ[71,3] = 12000
$blockHQ[71,4] = 12000
$blockHQ[71,5] = 12000
$blockHQ[71,6] = '$/caja 15 unidades'
$blockHQ[71,7] = 'Región Metropolitana'
$blockHQ[71,8] = 800
$blockHQ[71,9] = 15
$blockHQ[72,0] = 'Marina'
$blockHQ[72,1] = 'Primera'
$blockHQ[72,2] = 300
$blockHQ[72,3] = 9000
$blockHQ[72,4] = 9000
$blockHQ[72,5] = 9000
$blockHQ[72,6] = '$/caja 15 unidades'
$blockHQ[72,7] = 'Región Metropolitana'
$blockHQ[72,8] = 600
$blockHQ[72,9] = 15
$blockHQ[73,0] = 'Escarola'
$blockHQ[73,1] = 'Primera'
$blockHQ[73,2] = 400
$blockHQ[73,3] = 14000
$blockHQ[73,4] = 15000
$blockHQ[73,5] = 14500
$blockHQ[73,6] = '$/caja 15 unidades'
$blockHQ[73,7] = 'Región de Coquimbo'
$blockHQ[73,8] = 967
$blockHQ[73,9] = 15
$blockHQ[74,0] = 'Escarola'
$blockHQ[74,1] = 'Segunda'
$blockHQ[74,2] = 200
$blockHQ[74,3] = 12000
$blockHQ[74,4] = 12000
$blockHQ[74,5] = 12000
$blockHQ[74,6] = '$/caja 18 unidades'
$blockHQ[74,7] = 'Región de Coquimbo'
$blockHQ[74,8] = 667
$blockHQ[74,9] = 18
$blockHQ[75,0] = 'Escarola'
$blockHQ[75,1] = 'Primera'
$blockHQ[75,2] = 120
$blockHQ[75,3] = 12000
$blockHQ[75,4] = 12000
$blockHQ[75,5] = 12000
$blockHQ[75,6] = '$/caja 15 unidades'
$blockHQ[75,7] = 'Región de Coquimbo'
$blockHQ[75,8] = 800
$blockHQ[75,9] = 15
$blockHQ[76,0] = 'Conconina(o)'
$blockHQ[76,1] = 'Primera'
$blockHQ[76,2] = 90
$blockHQ[76,3] = 8000
$blockHQ[76,4] = 8000
$blockHQ[76,5] = 8000
$blockHQ[76,6] = '$/caja 15 unidades'
$blockHQ[76,7] = 'Región Metropolitana'
$blockHQ[76,8] = 533
$blockHQ[76,9] = 15
$blockHQ[77,0] = 'Escarola'
$blockHQ[77,1] = 'Primera'
$blockHQ[77,2] = 400
$blockHQ[77,3] = 8500
$blockHQ[77,4] = 8500
$blockHQ[77,5] = 8500
$blockHQ[77,6] = '$/caja 15 unidades'
$blockHQ[77,7] = 'Región del Maule'
$blockHQ[77,8] = 567
$blockHQ[77,9] = 15
$blockHQ[78,0] = 'Escarola'
$blockHQ[78,1] = 'Primera'
$blockHQ[78,2] = 200
$blockHQ[78,3] = 12000
$blockHQ[78,4] = 12000
$blockHQ[78,5] = 12000
$blockHQ[78,6] = '$/caja 15 unidades'
$blockHQ[78,7] = 'Región del Maule'
$blockHQ[78,8] = 800
$blockHQ[78,9] = 15
$blockHQ[79,0] = 'Escarola'
$blockHQ[79,1] = 'Primera'
$blockHQ[79,2] = 200
$blockHQ[79,3] = 14000
$blockHQ[79,4] = 14000
$blockHQ[79,5] = 14000
$blockHQ[79,6] = '$/caja 15 unidades'
$blockHQ[79,7] = 'Región de Coquimbo'
$blockHQ[79,8] = 933
$blockHQ[79,9] = 15
$blockHQ[80,0] = 'Escarola'
$blockHQ[80,1] = 'Primera'
$blockHQ[80,2] = 200
$blockHQ[80,3] = 10000
$blockHQ[80,4] = 10000
$blockHQ[80,5] = 10000
$blockHQ[80,6] = '$/caja 15 unidades'
$blockHQ[80,7] = 'Región de La Araucanía'
$blockHQ[80,8] = 667
$blockHQ[80,9] = 15
$blockHQ[81,0] = 'Escarola'
$blockHQ[81,1] = 'Primera'
$blockHQ[81,2] = 350
$blockHQ[81,3] = 12000
$blockHQ[81,4] = 12000
$blockHQ[81,5] = 12000
$blockHQ[81,6] = '$/caja 15 unidades'
$blockHQ[81,7] = 'Región de Coquimbo'
$blockHQ[81,8] = 800
$blockHQ[81,9] = 15
$blockHQ[82,0] = 'Marina'
$blockHQ[82,1] = 'Segunda'
$blockHQ[82,2] = 200
$blockHQ[82,3] = 9000
$blockHQ[82,4] = 9000
$blockHQ[82,5] = 9000
$blockHQ[82,6] = '$/caja 18 unidades'
$blockHQ[82,7] = 'Región Metropolitana'
$blockHQ[82,8] = 500
$blockHQ[82,9] = 18
$blockHQ[83,0] = 'Escarola'
$blockHQ[83,1] = 'Primera'
$blockHQ[83,2] = 200
$blockHQ[83,3] = 11000
$blockHQ[83,4] = 11000
$blockHQ[83,5] = 11000
$blockHQ[83,6] = '$/caja 15 unidades'
$blockHQ[83,7] = 'Región del Maule'
$blockHQ[83,8] = 733
$blockHQ[83,9] = 15
$blockHQ[84,0] = 'Escarola'
$blockHQ[84,1] = 'Primera'
$blockHQ[84,2] = 500
$blockHQ[84,3] = 7500
$blockHQ[84,4] = 8000
$blockHQ[84,5] = 7750
$blockHQ[84,6] = '$/caja 15 unidades'
$blockHQ[84,7] = 'Región Metropolitana'
$blockHQ[84,8] = 517
$blockHQ[84,9] = 15
$ws.Range("H935:Q1019").Value = $blockHQ

# Newly created rows 1018 and 1019 need the constant metadata columns filled in
$constABC = New-Object "object[,]" 2,3
$constABC[0,0] = 4
$constABC[0,1] = 'Feria Lagunitas de Puerto Montt'
$constABC[0,2] = 'Los Lagos'
$constABC[1,0] = 4
$constABC[1,1] = 'Feria Lagunitas de Puerto Montt'
$constABC[1,2] = 'Los Lagos'
$ws.Range("A1018:C1019").Value = $constABC

$constEFG = New-Object "object[,]" 2,3
$constEFG[0,0] = 10
$constEFG[0,1] = 100112033
$constEFG[0,2] = 'Lechuga'
$constEFG[1,0] = 10
$constEFG[1,1] = 100112033
$constEFG[1,2] = 'Lechuga'
$ws.Range("E1018:G1019").Value = $constEFG

$ws.Range("R1018:R1019").Value = 'Hortaliza'

$ws.Range("D1018:D1019").NumberFormat = "YYYY-MM-DD HH:MM:SS"

